$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2 = 'yale_blue'
    3 = 'oxford_blue'
    4 = 'charcoal'
    5 = 'sea_blue'
    6 = 'ucla_blue'
    7 = 'dark_pastel_blue'
    8 = 'dark_slate_gray'
    10 = 'medium_carmine'
    11 = 'coffee'
    12 = 'fuzzy_wuzzy'
    13 = 'eggplant'
    14 = 'caput_mortuum'
    15 = 'dark_midnight_blue'
    16 = 'oxford_blue'
    17 = 'charcoal'
    18 = 'sea_blue'
    20 = 'dark_pastel_blue'
    21 = 'msu_green'
    23 = 'burnt_umber'
    24 = 'coffee'
    26 = 'purple_taupe'
    27 = 'dark_sienna'
    29 = 'dark_pastel_blue'
    30 = 'ucla_blue'
    31 = 'ucla_blue'
    32 = 'sea_blue'
    33 = 'dark_midnight_blue'
    34 = 'purple_taupe'
    35 = 'teal_green'
    36 = 'dark_slate_gray'
    37 = 'dark_jungle_green'
    38 = 'khaki'
    39 = 'pale_taupe'
    40 = 'pale_brown'
    41 = 'café_au_lait'
    42 = 'medium_carmine'
    43 = 'coffee'
    44 = 'dark_lava'
    46 = 'deep_carmine'
    47 = 'rose_vale'
    48 = 'deep_coffee'
    49 = 'pale_silver'
    50 = 'davy_grey'
    51 = 'rose_quartz'
    52 = 'lavender_gray'
    53 = 'dark_gray'
    54 = 'lavender_gray'
    55 = 'onyx'
    56 = 'dim_gray'
    57 = 'shadow'
    58 = 'battleship_grey'
    59 = 'dark_tan'
    60 = 'dark_gray'
    61 = 'onyx'
    62 = 'dark_lava'
    63 = 'onyx'
    64 = 'onyx'
    65 = 'bistre'
    66 = 'onyx'
    67 = 'dark_lava'
    68 = 'onyx'
    69 = 'dim_gray'
    70 = 'dark_gray'
    71 = 'purple_taupe'
    72 = 'oxford_blue'
    73 = 'outer_space'
    74 = 'onyx'
    75 = 'onyx'
    76 = 'outer_space'
    77 = 'liver'
    78 = 'charcoal'
    79 = 'charcoal'
    80 = 'prussian_blue'
    81 = 'charcoal'
    82 = 'liver'
    83 = 'oxford_blue'
    85 = 'davy_grey'
    86 = 'oxford_blue'
    87 = 'purple_taupe'
    88 = 'dim_gray'
    89 = 'outer_space'
    90 = 'dark_cerulean'
    92 = 'davy_grey'
    93 = 'outer_space'
    94 = 'charcoal'
    95 = 'onyx'
    96 = 'khaki'
    97 = 'dark_lava'
    98 = 'bistre'
    100 = 'charcoal'
    101 = 'dim_gray'
    102 = 'ucla_blue'
    103 = 'dark_midnight_blue'
    104 = 'grullo'
    105 = 'shadow'
    107 = 'onyx'
    108 = 'onyx'
    109 = 'onyx'
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 14).Value = $values[$row]
}

$wb.Save()
